# Add a new, blank slide right after the existing slide.
#
# We duplicate slide 1 (whose <p:cSld> wrapper already carries the
# standard boilerplate PowerPoint emits for a freshly authored slide -
# a zeroed <p:grpSpPr><a:xfrm>, a <p:extLst> with a p14:creationId, and
# a <p:clrMapOvr><a:masterClrMapping/>) and then strip out the picture
# shape it inherited, leaving an otherwise empty slide that uses the
# same (Blank) layout as slide 1.

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

$newSlide = $s1.Duplicate()

$s2 = $p.Slides.Item(2)
while ($s2.Shapes.Count -gt 0) {
    $s2.Shapes.Item(1).Delete()
}
